$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C) for rows 2-15 from 45243 to 45244
for ($row = 2; $row -le 15; $row++) {
    $addr = "C" + $row
    $ws.Range($addr).Value = 45244
}
